$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The participant list in column A currently has 4 rows:
#   A1 emaild
#   A2 cindy.leschaud@gmail.com
#   A3 cindy@designpond.ch   (mailto: hyperlink, "Lien hypertexte" style)
#   A4 info@designpond.ch    (mailto: hyperlink, "Lien hypertexte" style)
#
# Remove the hyperlink on A3 before we overwrite its value.
$null = $ws.Range("A3").Hyperlinks.Delete()

# Clear every cell (content + the hyperlink formatting/style) so the
# now-unused "Lien hypertexte" style stops being referenced by any cell.
$null = $ws.Cells.Delete()

# Re-enter the participant emails: keep the first two, replace the third
# participant with a new one, and drop the fourth participant entirely.
$ws.Range("A1").Value = "emaild"
$ws.Range("A2").Value = "cindy.leschaud@gmail.com"
$ws.Range("A3").Value = "hello@yahoo.fr"

# The "Lien hypertexte" cell style is no longer used by any cell - remove it.
$null = $wb.Styles.Item("Lien hypertexte").Delete()

# Leave the selection on row 2, matching the saved view state.
$null = $ws.Rows(2).Select()
